# Refresh the crypto price/volume table with the latest crawl values.
# A couple of rows also swapped rank order (e.g. Chainlink/Polygon,
# Toncoin/LEO, InjectiveProtocol/Bittensor, Stellar/Maker, Monero/ApeXProtocol)
# and row 51 now reports "dogwifhat" instead of "Stacks".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.268.50"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "3.729.42"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "192.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.02%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.731"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "60.40"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +12.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000292"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("D14").Value = "4.335.25"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").Value = "3.733.80"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("D20").Value = "69.149.82"
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "413.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "90.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.51%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.47%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.123"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.81%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "640.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.03%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "45.73"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "67.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.79%  "
$ws.Range("D37").Value = "0.0₃0841"
$ws.Range("E37").Value = "  -9.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.417"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.43%  "
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.141"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.94%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.900.31"
$ws.Range("E45").Value = "  +4.61%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.140"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.81%  "
$ws.Range("E47").Value = "  -2.41%  "
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -17.25%  "
